$d = $word.ActiveDocument

# Replace the topic cell text
$d.Content.Find.Execute("Accessing & Using Databases", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Accessing & Using External Data", 2)

# Replace the tools cell text (exact match "SQL" only)
$d.Content.Find.Execute("SQL", $true, $false, $false, $false, $false,
                         $true, 1, $false, "SQL, other APIs", 2)
